$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from DICOM")

# --- Sheet "Metadata" updates ---

# Version 1.0.0 -> 1.1.0
$ws1.Range("B3").Value = "1.1.0"

# Experimental value now "false" (stored as text, not boolean)
$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "false"

# Date updated
$ws1.Range("B8").Value = "2024-06-20T08:51:57-05:00"

# Contact rows: row10 keeps "Contact" label but new value; rows 11-12 are new rows
$ws1.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/radiology/)"

$ws1.Range("A11").Value = "Contact"
$ws1.Range("B11").Value = "null (radiology@ihe.net)"
$ws1.Range("A11").Style = $ws1.Range("A10").Style
$ws1.Range("B11").Style = $ws1.Range("B10").Style

$ws1.Range("A12").Value = "Contact"
$ws1.Range("B12").Value = "IHE Radiology Technical Committee (radiology@ihe.net)"
$ws1.Range("A12").Style = $ws1.Range("A10").Style
$ws1.Range("B12").Style = $ws1.Range("B10").Style

# Jurisdiction row (shift content down by two rows due to inserted contact rows)
$ws1.Range("A13").Value = "Jurisdiction"
$ws1.Range("B13").Value = "Global (Whole world)"
$ws1.Range("A13").Style = $ws1.Range("A10").Style
$ws1.Range("B13").Style = $ws1.Range("B10").Style

# Description row
$ws1.Range("A14").Value = "Description"
$ws1.Range("B14").Value = "Codes representing the applicable endpoint payloadType to retrieve a study."
$ws1.Range("A14").Style = $ws1.Range("A10").Style
$ws1.Range("B14").Style = $ws1.Range("B10").Style

# Purpose row
$ws1.Range("A15").Value = "Purpose"
$ws1.Range("B15").Value = ""
$ws1.Range("A15").Style = $ws1.Range("A10").Style
$ws1.Range("B15").Style = $ws1.Range("B10").Style

# Copyright row
$ws1.Range("A16").Value = "Copyright"
$ws1.Range("B16").Value = ""
$ws1.Range("A16").Style = $ws1.Range("A10").Style
$ws1.Range("B16").Style = $ws1.Range("B10").Style

# Immutable row
$ws1.Range("A17").Value = "Immutable"
$ws1.Range("B17").Value = "BooleanType[null]"
$ws1.Range("A17").Style = $ws1.Range("A10").Style
$ws1.Range("B17").Style = $ws1.Range("B10").Style
